# Insert 3 new data rows (595, 596, 597) into the sheet, shifting the
# existing rows 595-662 down to 598-665. This matches the commit's weekly
# price update: a new reporting date (serial 45142, i.e. 2023-08-04) is
# prepended with 3 quality/variety records for "Alcachofa".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 595:597 -> 598:600 (and everything below accordingly) by
# inserting 3 blank rows at row 595.
$ws.Rows("595:597").Insert()

# New row 595: Española / Extra
$ws.Cells.Item(595, 1).Value = 9
$ws.Cells.Item(595, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(595, 3).Value = "Metropolitana"
$ws.Cells.Item(595, 4).Value = 45142
$ws.Cells.Item(595, 5).Value = 13
$ws.Cells.Item(595, 6).Value = 100112013
$ws.Cells.Item(595, 7).Value = "Alcachofa"
$ws.Cells.Item(595, 8).Value = "Española"
$ws.Cells.Item(595, 9).Value = "Extra"
$ws.Cells.Item(595, 10).Value = 52
$ws.Cells.Item(595, 11).Value = 15000
$ws.Cells.Item(595, 12).Value = 16000
$ws.Cells.Item(595, 13).Value = 15500
$ws.Cells.Item(595, 14).Value = "$/caja 25 unidades"
$ws.Cells.Item(595, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(595, 16).Value = 15500
$ws.Cells.Item(595, 17).Value = 1
$ws.Cells.Item(595, 18).Value = "Hortaliza"

# New row 596: Española / Primera
$ws.Cells.Item(596, 1).Value = 9
$ws.Cells.Item(596, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(596, 3).Value = "Metropolitana"
$ws.Cells.Item(596, 4).Value = 45142
$ws.Cells.Item(596, 5).Value = 13
$ws.Cells.Item(596, 6).Value = 100112013
$ws.Cells.Item(596, 7).Value = "Alcachofa"
$ws.Cells.Item(596, 8).Value = "Española"
$ws.Cells.Item(596, 9).Value = "Primera"
$ws.Cells.Item(596, 10).Value = 70
$ws.Cells.Item(596, 11).Value = 14000
$ws.Cells.Item(596, 12).Value = 15000
$ws.Cells.Item(596, 13).Value = 14500
$ws.Cells.Item(596, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(596, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(596, 16).Value = 483
$ws.Cells.Item(596, 17).Value = 30
$ws.Cells.Item(596, 18).Value = "Hortaliza"

# New row 597: Madrigal / Primera
$ws.Cells.Item(597, 1).Value = 9
$ws.Cells.Item(597, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(597, 3).Value = "Metropolitana"
$ws.Cells.Item(597, 4).Value = 45142
$ws.Cells.Item(597, 5).Value = 13
$ws.Cells.Item(597, 6).Value = 100112013
$ws.Cells.Item(597, 7).Value = "Alcachofa"
$ws.Cells.Item(597, 8).Value = "Madrigal"
$ws.Cells.Item(597, 9).Value = "Primera"
$ws.Cells.Item(597, 10).Value = 52
$ws.Cells.Item(597, 11).Value = 12000
$ws.Cells.Item(597, 12).Value = 13000
$ws.Cells.Item(597, 13).Value = 12500
$ws.Cells.Item(597, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(597, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(597, 16).Value = 312
$ws.Cells.Item(597, 17).Value = 40
$ws.Cells.Item(597, 18).Value = "Hortaliza"
